# Refresh the "cryptos" price list: Price (D) and Volume(1h) (E) columns are
# updated for most rows, and rows 45/46 (EnergySwap / PancakeSwap) swap
# their Coin/Link/Price data along with new volume figures.
#
# Price/volume figures are stored as plain text in the sheet (not numbers),
# so values that would otherwise parse as a number (e.g. "288.32") are
# entered with a leading apostrophe - Excel's normal "force text" quote
# prefix - to keep them as text instead of being auto-converted to Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '21.655.48'
$ws.Range('E2').Value = '  -1.65%  '

$ws.Range('D3').Value = '1.533.27'

$ws.Range('E4').Value = '  +0.15%  '

$ws.Range('D6').Value = '''288.32'
$ws.Range('E6').Value = '  +0.47%  '

$ws.Range('D7').Value = '''0.3946'
$ws.Range('E7').Value = '  +3.16%  '

$ws.Range('D8').Value = '''0.3154'
$ws.Range('E8').Value = '  -2.63%  '

$ws.Range('D9').Value = '''42.21'
$ws.Range('E9').Value = '  +2.02%  '

$ws.Range('D10').Value = '''0.07151'
$ws.Range('E10').Value = '  -2.26%  '

$ws.Range('D11').Value = '''1.045'
$ws.Range('E11').Value = '  -6.88%  '

$ws.Range('D12').Value = '''1.001'
$ws.Range('E12').Value = '  +0.17%  '

$ws.Range('D13').Value = '''5.624'
$ws.Range('E13').Value = '  -1.71%  '

$ws.Range('D14').Value = '''18.48'
$ws.Range('E14').Value = '  -4.62%  '

$ws.Range('D15').Value = '''6.606'

$ws.Range('D16').Value = '1.532.81'
$ws.Range('E16').Value = '  -1.38%  '

$ws.Range('D17').Value = '''0.00001092'
$ws.Range('E17').Value = '  -0.02%  '

$ws.Range('D18').Value = '''0.06592'
$ws.Range('E18').Value = '  -0.47%  '

$ws.Range('D19').Value = '''82.97'
$ws.Range('E19').Value = '  -2.65%  '

$ws.Range('E20').Value = '  +0.12%  '

$ws.Range('D21').Value = '''6.098'
$ws.Range('E21').Value = '  -4.82%  '

$ws.Range('D22').Value = '''15.39'
$ws.Range('E22').Value = '  -3.50%  '

$ws.Range('D23').Value = '''10.82'
$ws.Range('E23').Value = '  -5.52%  '

$ws.Range('E24').Value = '  +4.14%  '

$ws.Range('D25').Value = '21.657.31'
$ws.Range('E25').Value = '  -1.70%  '

$ws.Range('D26').Value = '''2.332'
$ws.Range('E26').Value = '  -7.72%  '

$ws.Range('D27').Value = '''147.64'
$ws.Range('E27').Value = '  -0.85%  '

$ws.Range('D28').Value = '''18.27'
$ws.Range('E28').Value = '  -2.90%  '

$ws.Range('D29').Value = '''4.832'
$ws.Range('E29').Value = '  -0.48%  '

$ws.Range('D30').Value = '1.716.78'
$ws.Range('E30').Value = '  -0.65%  '

$ws.Range('D31').Value = '''116.96'
$ws.Range('E31').Value = '  -3.10%  '

$ws.Range('D32').Value = '''5.849'
$ws.Range('E32').Value = '  -0.67%  '

$ws.Range('D33').Value = '''0.9461'
$ws.Range('E33').Value = '  -13.77%  '

$ws.Range('D34').Value = '''0.08126'
$ws.Range('E34').Value = '  -0.11%  '

$ws.Range('D35').Value = '''8.505'

$ws.Range('D36').Value = '''0.06061'
$ws.Range('E36').Value = '  -2.36%  '

$ws.Range('D37').Value = '''5.088'
$ws.Range('E37').Value = '  -3.18%  '

$ws.Range('D38').Value = '''0.02201'
$ws.Range('E38').Value = '  -4.14%  '

$ws.Range('D39').Value = '''1.436'
$ws.Range('E39').Value = '  -12.56%  '

$ws.Range('D40').Value = '''0.2014'
$ws.Range('E40').Value = '  -4.19%  '

$ws.Range('D41').Value = '''1.177'
$ws.Range('E41').Value = '  -3.59%  '

$ws.Range('D42').Value = '''1.000'

$ws.Range('D43').Value = '''10.78'
$ws.Range('E43').Value = '  -0.70%  '

$ws.Range('D44').Value = '''0.5730'
$ws.Range('E44').Value = '  -3.40%  '

$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D45').Value = '''3.720'
$ws.Range('E45').Value = '  -0.02%  '

$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '''12.89'
$ws.Range('E46').Value = '  -4.33%  '

$ws.Range('D47').Value = '''0.5472'
$ws.Range('E47').Value = '  -4.77%  '

$ws.Range('D48').Value = '''1.152'
$ws.Range('E48').Value = '  -0.30%  '

$ws.Range('D49').Value = '''1.859'
$ws.Range('E49').Value = '  -3.81%  '

$ws.Range('D50').Value = '''115.46'
$ws.Range('E50').Value = '  -3.31%  '

$ws.Range('D51').Value = '''0.06687'
$ws.Range('E51').Value = '  -2.70%  '
